$wb = $excel.ActiveWorkbook

# --- Step 1: rename the existing "2022-Q1" sheet to "2022-Q4" ---
$sheetQ4 = $wb.Worksheets.Item("2022-Q1")
$sheetQ4.Name = "2022-Q4"

# --- Step 2: add a brand-new sheet right after it, and copy the old ---
# --- (now renamed) Q1 data into it, then rename it back to "2022-Q1" ---
$sheetQ1 = $wb.Worksheets.Add($null, $sheetQ4)
$sheetQ4.UsedRange.Copy($sheetQ1.Range("A1"))
$sheetQ1.Name = "2022-Q1"

# --- Step 3: overwrite "2022-Q4" sheet with the new quarterly data ---
$sheetQ4.Cells.Clear()

$sheetQ4.Range("B1").Value = "基金代码"
$sheetQ4.Range("C1").Value = "基金名称"
$sheetQ4.Range("D1").Value = "基金规模"
$sheetQ4.Range("E1").Value = "股票总仓位"
$sheetQ4.Range("F1").Value = "仓位占比"
$sheetQ4.Range("G1").Value = "持有市值(亿元)"
$sheetQ4.Range("H1").Value = "仓位排名"
$sheetQ4.Range("B1:H1").Font.Bold = $true

$q4data = @(
    @(0, "009686", "华夏磐利一年定期开放混合A", "11.49", "92.56", "5.08", "0.5837", 2),
    @(1, "015697", "华夏磐润两年定开混合A",     "2.68",  "86.76", "5.07", "0.1359", 1),
    @(2, "015698", "华夏磐润两年定开混合C",     "0.99",  "86.76", "5.07", "0.0502", 1),
    @(3, "009687", "华夏磐利一年定期开放混合C", "0.46",  "92.56", "5.08", "0.0234", 2),
    @(4, "164401", "前海开源中证健康产业指数",   "1.89",  "94.43", "1.12", "0.0212", 4)
)

$r = 2
foreach ($row in $q4data) {
    $sheetQ4.Cells.Item($r, 1).Value = $row[0]
    $sheetQ4.Cells.Item($r, 2).Value = "'" + $row[1]
    $sheetQ4.Cells.Item($r, 3).Value = "'" + $row[2]
    $sheetQ4.Cells.Item($r, 4).Value = "'" + $row[3]
    $sheetQ4.Cells.Item($r, 5).Value = "'" + $row[4]
    $sheetQ4.Cells.Item($r, 6).Value = "'" + $row[5]
    $sheetQ4.Cells.Item($r, 7).Value = "'" + $row[6]
    $sheetQ4.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# --- Step 4: update the "总计" (totals) sheet: insert a new row for ---
# --- 2022-Q4 above the existing 2022-Q1 row ---
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.8100000000000001
$total.Range("A3").Value = 1
